$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.006.33"
$ws.Range("E2").Value = "  -0.62%  "

$ws.Range("D3").Value = "1.584.55"
$ws.Range("E3").Value = "  -1.24%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.47%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.003"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.18%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "300.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3760"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3576"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.42%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "50.51"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +4.02%  "

$ws.Range("E10").Value = "  +0.48%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08013"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.214"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -4.83%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.87"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.452"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.54%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.295"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -5.06%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001219"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -3.74%  "

$ws.Range("D17").Value = "1.587.71"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.12"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.74%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06791"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.88"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -2.78%  "

$ws.Range("E21").Value = "  +0.24%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.420"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -2.49%  "

$ws.Range("E23").Value = "  -1.74%  "

$ws.Range("D24").Value = "23.006.83"
$ws.Range("E24").Value = "  -0.59%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.365"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.754"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.28%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.74"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.85%  "

$ws.Range("E28").Value = "  -1.99%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.197"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -1.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.56"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.345"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.66%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.493"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.63%  "

$ws.Range("D33").Value = "1.764.76"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9365"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.42%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07336"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -5.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02671"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.06%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.987"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.34%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.08751"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -1.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.050"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2470"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.335"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6867"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -3.95%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.85"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -7.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.91"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -5.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6384"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.36%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.984"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.241"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.59"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07879"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.189"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.185"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.62%  "

